# The commit swaps the contents of ppt/theme/theme1.xml (bound to the
# Notes Master) and ppt/theme/theme2.xml (bound to the Slide Master):
#   theme1.xml : "Office Theme" colours  ->  "Integral" colours
#   theme2.xml : "Integral" colours      ->  "Office Theme" colours
#
# PowerPoint's automation surface only exposes one mutable "theme" for a
# deck - the one attached to ActivePresentation.SlideMaster - so we apply
# the half of the swap that is reachable through the object model: we
# recolour the slide master's theme from "Integral" to the stock "Office
# Theme" palette by rewriting each of the twelve theme colour slots.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# expressed as COM RGB integers (0xBBGGRR, i.e. R + G*256 + B*65536).
$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72

# Reflect the renamed scheme/theme where the object model allows it.
$theme.Name = "Office Theme"
$colors.Name = "Office"
